$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.169.84'
$ws.Range("E2").Value = '  +4.57%  '

$ws.Range("D3").Value = '2.512.57'
$ws.Range("E3").Value = '  +2.98%  '

$ws.Range("E4").Value = '  -0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '495.37'
$ws.Range("E5").Value = '  +3.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.26'
$ws.Range("E6").Value = '  +12.12%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").Value = '  +3.27%  '

$ws.Range("D9").Value = '2.531.30'
$ws.Range("E9").Value = '  +3.35%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("E10").Value = '  +3.86%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.77'
$ws.Range("E11").Value = '  +6.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.338'
$ws.Range("E12").Value = '  +4.36%  '

$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("D14").Value = '2.946.77'
$ws.Range("E14").Value = '  +3.11%  '

$ws.Range("D15").Value = '57.285.64'
$ws.Range("E15").Value = '  +4.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.42'
$ws.Range("E16").Value = '  +4.79%  '

$ws.Range("E17").Value = '  +3.10%  '

$ws.Range("D18").Value = '2.528.56'
$ws.Range("E18").Value = '  +3.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.57'
$ws.Range("E19").Value = '  +5.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.36'
$ws.Range("E20").Value = '  +5.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.15'
$ws.Range("E21").Value = '  +3.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.94'
$ws.Range("E23").Value = '  +5.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.51'
$ws.Range("E24").Value = '  +2.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.412'
$ws.Range("E25").Value = '  +2.01%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.41%  '

$ws.Range("D28").Value = '2.613.74'
$ws.Range("E28").Value = '  +2.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.65'
$ws.Range("E29").Value = '  +5.08%  '

$ws.Range("D30").Value = '0.0₃0829'
$ws.Range("E30").Value = '  +7.08%  '

$ws.Range("E31").Value = '  +0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.60'
$ws.Range("E32").Value = '  +2.49%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.39'
$ws.Range("E33").Value = '  +3.00%  '

$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").Value = '  +4.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("E35").Value = '  +3.80%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.17'
$ws.Range("E36").Value = '  +5.60%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.85'
$ws.Range("E37").Value = '  +6.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.889'
$ws.Range("E38").Value = '  +4.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.42'
$ws.Range("E39").Value = '  +9.69%  '

$ws.Range("E40").Value = '  +3.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.55'
$ws.Range("E41").Value = '  +4.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.621'
$ws.Range("E42").Value = '  +4.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0565'
$ws.Range("E43").Value = '  +4.14%  '

$ws.Range("E44").Value = '  -0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.96'
$ws.Range("E45").Value = '  +6.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '268.91'
$ws.Range("E46").Value = '  +5.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0943'
$ws.Range("E47").Value = '  +4.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0231'
$ws.Range("E48").Value = '  +4.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.19'
$ws.Range("E49").Value = '  +0.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.12'
$ws.Range("E50").Value = '  +6.04%  '

$ws.Range("D51").Value = '1.900.59'
$ws.Range("E51").Value = '  -1.92%  '
